{"js": "const replacements = [\n  [\"2024-02-25 Sunday\", \"2024-02-26 Monday\"],\n  [\"202\u00d72=404\", \"400\u00d75=2000\"],\n  [\"930\u00d75=4650\", \"882\u00d72=1764\"],\n  [\"448\u00d73=1344\", \"452\u00d77=3164\"],\n  [\"923\u00d76=5538\", \"311\u00d79=2799\"],\n  [\"449\u00d79=4041\", \"355\u00d74=1420\"],\n  [\"509\u00d78=4072\", \"256\u00d72=512\"],\n  [\"704\u00d75=3520\", \"932\u00d76=5592\"],\n  [\"874\u00d77=6118\", \"325\u00d77=2275\"],\n  [\"171\u00d79=1539\", \"792\u00d77=5544\"],\n  [\"819\u00d77=5733\", \"190\u00d76=1140\"],\n  [\"151\u00d78=1208\", \"452\u00d75=2260\"],\n  [\"937\u00d72=1874\", \"532\u00d78=4256\"],\n  [\"512\u00d73=1536\", \"196\u00d73=588\"],\n  [\"790\u00d78=6320\", \"157\u00d76=942\"],\n  [\"622\u00d73=1866\", \"533\u00d75=2665\"],\n  [\"406\u00d75=2030\", \"218\u00d79=1962\"],\n  [\"177\u00d74=708\", \"799\u00d72=1598\"],\n  [\"224\u00d77=1568\", \"982\u00d73=2946\"],\n  [\"123\u00d77=861\", \"274\u00d72=548\"],\n  [\"337\u00d75=1685\", \"717\u00d76=4302\"],\n  [\"537\u00d77=3759\", \"119\u00d76=714\"],\n  [\"236\u00d73=708\", \"223\u00d75=1115\"],\n  [\"264\u00d72=528\", \"633\u00d79=5697\"],\n  [\"999\u00d74=3996\", \"264\u00d74=1056\"],\n  [\"667\u00d79=6003\", \"477\u00d74=1908\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-25 Sunday\", \"2024-02-26 Monday\"),\n    @(\"202\u00d72=404\", \"400\u00d75=2000\"),\n    @(\"930\u00d75=4650\", \"882\u00d72=1764\"),\n    @(\"448\u00d73=1344\", \"452\u00d77=3164\"),\n    @(\"923\u00d76=5538\", \"311\u00d79=2799\"),\n    @(\"449\u00d79=4041\", \"355\u00d74=1420\"),\n    @(\"509\u00d78=4072\", \"256\u00d72=512\"),\n    @(\"704\u00d75=3520\", \"932\u00d76=5592\"),\n    @(\"874\u00d77=6118\", \"325\u00d77=2275\"),\n    @(\"171\u00d79=1539\", \"792\u00d77=5544\"),\n    @(\"819\u00d77=5733\", \"190\u00d76=1140\"),\n    @(\"151\u00d78=1208\", \"452\u00d75=2260\"),\n    @(\"937\u00d72=1874\", \"532\u00d78=4256\"),\n    @(\"512\u00d73=1536\", \"196\u00d73=588\"),\n    @(\"790\u00d78=6320\", \"157\u00d76=942\"),\n    @(\"622\u00d73=1866\", \"533\u00d75=2665\"),\n    @(\"406\u00d75=2030\", \"218\u00d79=1962\"),\n    @(\"177\u00d74=708\", \"799\u00d72=1598\"),\n    @(\"224\u00d77=1568\", \"982\u00d73=2946\"),\n    @(\"123\u00d77=861\", \"274\u00d72=548\"),\n    @(\"337\u00d75=1685\", \"717\u00d76=4302\"),\n    @(\"537\u00d77=3759\", \"119\u00d76=714\"),\n    @(\"236\u00d73=708\", \"223\u00d75=1115\"),\n    @(\"264\u00d72=528\", \"633\u00d79=5697\"),\n    @(\"999\u00d74=3996\", \"264\u00d74=1056\"),\n    @(\"667\u00d79=6003\", \"477\u00d74=1908\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
